# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (interested count) column F on the "展览" and
# "全部类型" sheets to reflect newly generated numbers.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (sheet1) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3147
$ws1.Range("F4").Value = 1103
$ws1.Range("F5").Value = 90
$ws1.Range("F6").Value = 41
$ws1.Range("F8").Value = 41
$ws1.Range("F9").Value = 1133
$ws1.Range("F10").Value = 15859
$ws1.Range("F11").Value = 250
$ws1.Range("F14").Value = 6216
$ws1.Range("F16").Value = 113
$ws1.Range("F18").Value = 10
$ws1.Range("F22").Value = 635
$ws1.Range("F24").Value = 14
$ws1.Range("F25").Value = 8
$ws1.Range("F26").Value = 213
$ws1.Range("F27").Value = 873
$ws1.Range("F29").Value = 5009
$ws1.Range("F31").Value = 11118
$ws1.Range("F32").Value = 1236
$ws1.Range("F33").Value = 17
$ws1.Range("F34").Value = 131
$ws1.Range("F35").Value = 183

# ---- Sheet "全部类型" (sheet4) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3147
$ws4.Range("F5").Value = 1103
$ws4.Range("F6").Value = 90
$ws4.Range("F7").Value = 41
$ws4.Range("F9").Value = 41
$ws4.Range("F10").Value = 1133
$ws4.Range("F11").Value = 15859
$ws4.Range("F12").Value = 250
$ws4.Range("F15").Value = 6216
$ws4.Range("F17").Value = 113
$ws4.Range("F19").Value = 10
$ws4.Range("F23").Value = 635
$ws4.Range("F25").Value = 14
$ws4.Range("F26").Value = 8
$ws4.Range("F27").Value = 213
$ws4.Range("F28").Value = 873
$ws4.Range("F30").Value = 5009
$ws4.Range("F33").Value = 11118
$ws4.Range("F34").Value = 1236
$ws4.Range("F35").Value = 17
$ws4.Range("F36").Value = 131
$ws4.Range("F37").Value = 183
